# City index update: adds "okt" (October) 2025 monthly figures (column N)
# to the point-index sheet, refreshes the year-to-date city index summary
# row, and appends the next rolling 12-month window row to the sliding
# index sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "punktindeks_maned": add column N (October) values for the 2025
# rows that have data for that month.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("punktindeks_maned")

$octoberValues = @{
    5  = 0.33
    7  = 0.65
    9  = -9.82
    11 = 1.97
    13 = -0.94
    17 = 9.45
    19 = 1.14
    21 = -1.34
    23 = -1.13
    25 = -7.69
    27 = -3.55
    29 = 0.54
    31 = 5.78
    35 = 5.28
    37 = -6.27
    41 = -10.88
    43 = -5.73
    47 = 5.74
    55 = 2.68
    57 = -2.27
    59 = -0.96
}

foreach ($row in $octoberValues.Keys) {
    $ws2.Range("N$row").Value = $octoberValues[$row]
}

# ---------------------------------------------------------------------
# Sheet "byindeks_aarlig": refresh the rolling-year stats (row 2) and
# move the year-to-date summary (row 3) from "jan-sep" to "jan-okt".
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("byindeks_aarlig")

$ws3.Range("G2").Value = 0.6112575299850357

$ws3.Range("C3").Value = 10
$ws3.Range("D3").Value = -0.7
$ws3.Range("E3").Value = 0.9933
$ws3.Range("G3").Value = 0.875777907582678
$ws3.Range("K3").Value = "okt"
$ws3.Range("L3").Value = "jan-okt"
$ws3.Range("M3").Value = -2.4
$ws3.Range("N3").Value = 1

# ---------------------------------------------------------------------
# Sheet "by_glid_indeks": recompute the last rolling-12-month row (11)
# and append the new rolling 12-month window ending October 2025 (row 12).
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("by_glid_indeks")

$ws4.Range("A11").Value = 0.9919339311536522
$ws4.Range("B11").Value = -0.8066068846347751
$ws4.Range("E11").Value = 3.688100878101141
$ws4.Range("F11").Value = 0.9812146623273963
$ws4.Range("G11").Value = -2.8

$ws4.Range("A12").Value = 0.9910240340453687
$ws4.Range("B12").Value = -0.8975965954631349
$ws4.Range("C12").Value = 21
$ws4.Range("D12").Value = 13.64592445028583
$ws4.Range("E12").Value = 3.983051866642613
$ws4.Range("F12").Value = 1.078237620234392
$ws4.Range("G12").Value = -3.1
$ws4.Range("H12").Value = 1.4
$ws4.Range("I12").Value = "2023 - (nov 2024 - okt 2025)"
$ws4.Range("J12").Value = 45931
$ws4.Range("K12").Value = 10
$ws4.Range("L12").Value = 2025
$ws4.Range("M12").Value = "12_months"
